$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 260, shifting existing rows 260-347 down to 261-348.
$ws.Rows("260").Insert()

# Populate the newly inserted row 260 with the new weekly data point.
$ws.Range("A260").Value = 6
$ws.Range("B260").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C260").Value = "Metropolitana"
$ws.Range("D260").Value = 44559
$ws.Range("E260").Value = 13
$ws.Range("F260").Value = 100112039
$ws.Range("G260").Value = "Ciboulette"
$ws.Range("H260").Value = "Sin especificar"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 970
$ws.Range("K260").Value = 900
$ws.Range("L260").Value = 1000
$ws.Range("M260").Value = 945
$ws.Range("N260").Value = "`$/docena de atados"
$ws.Range("O260").Value = "Región Metropolitana"
$ws.Range("P260").Value = 315
$ws.Range("Q260").Value = 3
$ws.Range("R260").Value = "Hortaliza"
